$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to DAMSLTag (column I) and DialogAct (column J) following
# re-running SGNN dialog act annotation after transcript cleanup.

$ws.Range("I2").Value = "sd"
$ws.Range("J2").Value = "Statement-non-opinion"

$ws.Range("I9").Value = "sd"
$ws.Range("J9").Value = "Statement-non-opinion"

$ws.Range("I37").Value = "ba"
$ws.Range("J37").Value = "Appreciation"

$ws.Range("I39").Value = "sv"
$ws.Range("J39").Value = "Statement-opinion"

$ws.Range("I45").Value = "sd"
$ws.Range("J45").Value = "Statement-non-opinion"

$ws.Range("I50").Value = "ba"
$ws.Range("J50").Value = "Appreciation"

$ws.Range("I56").Value = "sd"
$ws.Range("J56").Value = "Statement-non-opinion"

$ws.Range("I61").Value = "ba"
$ws.Range("J61").Value = "Appreciation"
